$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15251.214
$ws.Range("H23").Value = 15251.214
$ws.Range("H51").Value = 2786
$ws.Range("J51").Value = 2942.8572
$ws.Range("L51").Value = 2942.8572
$ws.Range("N51").Value = -3910.8572
$ws.Range("H98").Value = 280523.44
$ws.Range("I98").Value = 329307
$ws.Range("J98").Value = 4083.3333
$ws.Range("K98").Value = 329307
$ws.Range("L98").Value = 4083.3333
$ws.Range("M98").Value = -327809
$ws.Range("N98").Value = -7079.3333
$ws.Range("H112").Value = 7576893
$ws.Range("I112").Value = 500
$ws.Range("K112").Value = 1500
$ws.Range("M112").Value = -392
$ws.Range("H122").Value = 280523.44
$ws.Range("I122").Value = 329307
$ws.Range("J122").Value = 4083.3333
$ws.Range("K122").Value = 987921
$ws.Range("L122").Value = 12249.9999
$ws.Range("M122").Value = -985471
$ws.Range("N122").Value = -17149.9999
$ws.Range("H137").Value = 1169.125
$ws.Range("I137").Value = 680.4231
$ws.Range("K137").Value = 2041.2693
$ws.Range("M137").Value = 508.7307000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4499
$ws.Range("I32").Value = 3140.258
$ws.Range("J32").Value = 6976.706
$ws.Range("K32").Value = 3140.258
$ws.Range("L32").Value = 6976.706
$ws.Range("M32").Value = -2853.258
$ws.Range("N32").Value = -7550.706
$ws.Range("H61").Value = 1092.1072
$ws.Range("I61").Value = 857.6957
$ws.Range("J61").Value = 2170.4
$ws.Range("K61").Value = 857.6957
$ws.Range("L61").Value = 2170.4
$ws.Range("M61").Value = -645.6957
$ws.Range("N61").Value = -2594.4
$ws.Range("H132").Value = 2303.4138
$ws.Range("I132").Value = 1993.1923
$ws.Range("J132").Value = 4992
$ws.Range("K132").Value = 5979.5769
$ws.Range("L132").Value = 14976
$ws.Range("M132").Value = -3449.5769
$ws.Range("N132").Value = -20036
$ws.Range("H136").Value = 1092.1072
$ws.Range("I136").Value = 857.6957
$ws.Range("J136").Value = 2170.4
$ws.Range("K136").Value = 2573.0871
$ws.Range("L136").Value = 6511.200000000001
$ws.Range("M136").Value = -23.08709999999974
$ws.Range("N136").Value = -11611.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 15007.875
$ws.Range("I19").Value = 15007.875
$ws.Range("K19").Value = 15007.875
$ws.Range("M19").Value = -14834.875
$ws.Range("H105").Value = 16669502
$ws.Range("I105").Value = 19610560
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 19610560
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -19608813
$ws.Range("N105").Value = -6994
$ws.Range("H137").Value = 51500
$ws.Range("J137").Value = 51500
$ws.Range("L137").Value = 51500
$ws.Range("N137").Value = -61700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5209648.5
$ws.Range("I99").Value = 10417765
$ws.Range("J99").Value = 1531.8334
$ws.Range("K99").Value = 10417765
$ws.Range("L99").Value = 1531.8334
$ws.Range("M99").Value = -10416267
$ws.Range("N99").Value = -4527.8334
$ws.Range("H126").Value = 5209648.5
$ws.Range("I126").Value = 10417765
$ws.Range("J126").Value = 1531.8334
$ws.Range("K126").Value = 31253295
$ws.Range("L126").Value = 4595.5002
$ws.Range("M126").Value = -31250825
$ws.Range("N126").Value = -9535.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 789.57776
$ws.Range("I107").Value = 151.19048
$ws.Range("J107").Value = 1348.1666
$ws.Range("K107").Value = 453.5714400000001
$ws.Range("L107").Value = 4044.4998
$ws.Range("M107").Value = 1466.42856
$ws.Range("N107").Value = -7884.4998
$ws.Range("H117").Value = 790
$ws.Range("I117").Value = 387.5
$ws.Range("J117").Value = 2400
$ws.Range("K117").Value = 1162.5
$ws.Range("L117").Value = 7200
$ws.Range("M117").Value = 2279.5
$ws.Range("N117").Value = -14084
$ws.Range("H122").Value = 720.8570999999999
$ws.Range("I122").Value = 535.85
$ws.Range("J122").Value = 1183.375
$ws.Range("K122").Value = 4822.650000000001
$ws.Range("L122").Value = 10650.375
$ws.Range("M122").Value = -2372.650000000001
$ws.Range("N122").Value = -15550.375
$ws.Range("H131").Value = 2508.175
$ws.Range("I131").Value = 519.75
$ws.Range("J131").Value = 2729.111
$ws.Range("K131").Value = 1559.25
$ws.Range("L131").Value = 8187.333
$ws.Range("M131").Value = 3480.75
$ws.Range("N131").Value = -18267.333
$ws.Range("H137").Value = 4814339
$ws.Range("I137").Value = 7146970
$ws.Range("J137").Value = 149076.14
$ws.Range("K137").Value = 21440910
$ws.Range("L137").Value = 447228.42
$ws.Range("M137").Value = -21435810
$ws.Range("N137").Value = -457428.42
$ws.Range("H138").Value = 1014
$ws.Range("I138").Value = 956.6667
$ws.Range("J138").Value = 1530
$ws.Range("K138").Value = 2870.0001
$ws.Range("L138").Value = 4590
$ws.Range("M138").Value = 2269.9999
$ws.Range("N138").Value = -14870
$ws.Range("H139").Value = 1687.619
$ws.Range("I139").Value = 1687.619
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5062.857
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 77.14300000000003
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 250501.5
$ws.Range("I12").Value = 250501.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 250501.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -250361.5
$ws.Range("N12").ClearContents()
$ws.Range("H102").Value = 1229.3903
$ws.Range("I102").Value = 1101.5625
$ws.Range("J102").Value = 1683.8889
$ws.Range("K102").Value = 1101.5625
$ws.Range("L102").Value = 1683.8889
$ws.Range("M102").Value = 520.4375
$ws.Range("N102").Value = -4927.8889
$ws.Range("H113").Value = 1520.0952
$ws.Range("I113").Value = 1548
$ws.Range("J113").Value = 1508.9333
$ws.Range("K113").Value = 1548
$ws.Range("L113").Value = 1508.9333
$ws.Range("M113").Value = 622
$ws.Range("N113").Value = -5848.9333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 45002
$ws.Range("J25").Value = 30008
$ws.Range("L25").Value = 30008
$ws.Range("N25").Value = -30468
$ws.Range("H136").Value = 3796.394
$ws.Range("I136").Value = 1232
$ws.Range("J136").Value = 8925.182000000001
$ws.Range("K136").Value = 3696
$ws.Range("L136").Value = 26775.546
$ws.Range("M136").Value = -1146
$ws.Range("N136").Value = -31875.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11495499
$ws.Range("I136").Value = 16667306
$ws.Range("K136").Value = 50001918
$ws.Range("M136").Value = -49999368
